# Auto-generated PowerShell COM-interop script
# Applies the cell-value updates described by the target diff
# (Sheets/Valefor_Profits.xlsx -> per-leve market price / profit recompute).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$updates = @{
    "H4" = 400
    "I4" = 400
    "J4" = 400
    "K4" = 400
    "L4" = 400
    "M4" = -286
    "N4" = -628
    "H18" = 1317.8462
    "I18" = 1239.5454
    "K18" = 1239.5454
    "M18" = -955.5454
    "H74" = 5740.7393
    "I74" = 6703.636
    "J74" = 4858.0835
    "K74" = 6703.636
    "L74" = 4858.0835
    "M74" = -5767.636
    "N74" = -6730.0835
    "H77" = 5740.7393
    "I77" = 6703.636
    "J77" = 4858.0835
    "K77" = 33518.18
    "L77" = 24290.4175
    "M77" = -28838.18
    "N77" = -33650.4175
    "H103" = 585423.7
    "I103" = 1852491.6
    "J103" = 623.0769
    "K103" = 5557474.800000001
    "L103" = 1869.2307
    "M103" = -5556888.800000001
    "N103" = -3041.2307
    "H127" = 2044.8422
    "J127" = 2160.1333
    "L127" = 6480.3999
    "N127" = -16400.3999
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$updates = @{
    "H109" = 29800
    "J109" = 29800
    "L109" = 29800
    "N109" = -32574
    "H112" = 10434.25
    "J112" = 10434.25
    "L112" = 10434.25
    "N112" = -13388.25
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$updates = @{
    "H110" = 40000
    "J110" = 40000
    "L110" = 40000
    "N110" = -48180
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$updates = @{
    "H22" = 23809702
    "I22" = 33333414
    "J22" = 417.33334
    "K22" = 33333414
    "L22" = 417.33334
    "M22" = -33333064
    "N22" = -1117.33334
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$updates = @{
    "H2" = 5440.5264
    "I2" = 12.25
    "J2" = 9388.362999999999
    "K2" = 73.5
    "L2" = 56330.178
    "M2" = 39.5
    "N2" = -56556.178
    "H15" = 267.5
    "I15" = 190
    "J15" = 500
    "K15" = 570
    "L15" = 1500
    "M15" = -430
    "N15" = -1780
    "H16" = 7783.3335
    "I16" = 1000
    "J16" = 9140
    "K16" = 3000
    "L16" = 27420
    "M16" = -2827
    "N16" = -27766
    "H17" = 631.25
    "I17" = 512.5
    "J17" = 750
    "K17" = 1537.5
    "L17" = 2250
    "M17" = -1368.5
    "N17" = -2588
    "H34" = 1421.1666
    "I34" = 125.5
    "J34" = 1680.3
    "K34" = 376.5
    "L34" = 5040.9
    "M34" = -292.5
    "N34" = -5208.9
    "H40" = 324.0909
    "I40" = 118.111115
    "J40" = 1251
    "K40" = 472.44446
    "L40" = 5004
    "M40" = -403.44446
    "N40" = -5142
    "H42" = 1251
    "I42" = 1000
    "J42" = 2004
    "K42" = 3000
    "L42" = 6012
    "M42" = -2466
    "N42" = -7080
    "H43" = 5950
    "J43" = 5950
    "L43" = 17850
    "N43" = -18078
    "H44" = 825.6667
    "I44" = 250
    "J44" = 1401.3334
    "K44" = 750
    "L44" = 4204.0002
    "M44" = -352
    "N44" = -5000.0002
    "H46" = 3000
    "I46" = 0
    "J46" = 3000
    "K46" = 0
    "L46" = 9000
    "N46" = -9182
    "H49" = 3600
    "J49" = 4000
    "L49" = 12000
    "N49" = -12312
    "H62" = 3897.75
    "I62" = 2000
    "J62" = 4530.3335
    "K62" = 6000
    "L62" = 13591.0005
    "M62" = -5314
    "N62" = -14963.0005
    "H63" = 4326.8335
    "I63" = 2659.6667
    "J63" = 5994
    "K63" = 7979.000100000001
    "L63" = 17982
    "M63" = -7230.000100000001
    "N63" = -19480
    "H64" = 3737.375
    "I64" = 499.5
    "J64" = 4816.6665
    "K64" = 1498.5
    "L64" = 14449.9995
    "M64" = -1228.5
    "N64" = -14989.9995
    "H65" = 3897.75
    "I65" = 2000
    "J65" = 4530.3335
    "K65" = 18000
    "L65" = 40773.0015
    "M65" = -14568
    "N65" = -47637.0015
    "H66" = 4326.8335
    "I66" = 2659.6667
    "J66" = 5994
    "K66" = 23937.0003
    "L66" = 53946
    "M66" = -20193.0003
    "N66" = -61434
    "H67" = 3737.375
    "I67" = 499.5
    "J67" = 4816.6665
    "K67" = 1498.5
    "L67" = 14449.9995
    "M67" = -562.5
    "N67" = -16321.9995
    "H69" = 1090
    "I69" = 872.2857
    "J69" = 1598
    "K69" = 2616.8571
    "L69" = 4794
    "M69" = -1805.8571
    "N69" = -6416
    "H70" = 2754
    "I70" = 302.8
    "J70" = 4115.778
    "K70" = 908.4000000000001
    "L70" = 12347.334
    "M70" = -593.4000000000001
    "N70" = -12977.334
    "H72" = 1090
    "I72" = 872.2857
    "J72" = 1598
    "K72" = 7850.571300000001
    "L72" = 14382
    "M72" = -3794.571300000001
    "N72" = -22494
    "H73" = 2754
    "I73" = 302.8
    "J73" = 4115.778
    "K73" = 908.4000000000001
    "L73" = 12347.334
    "M73" = 183.5999999999999
    "N73" = -14531.334
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
$clears = @("M46")
foreach ($cellRef in $clears) {
    $ws.Range($cellRef).ClearContents()
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$updates = @{
    "H102" = 1796.5476
    "I102" = 1160.3846
    "J102" = 2830.3125
    "K102" = 1160.3846
    "L102" = 2830.3125
    "M102" = 461.6153999999999
    "N102" = -6074.3125
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$updates = @{
    "H22" = 372.8857
    "I22" = 292.6316
    "J22" = 468.1875
    "K22" = 292.6316
    "L22" = 468.1875
    "M22" = 2.368400000000008
    "N22" = -1058.1875
    "H27" = 372.8857
    "I27" = 292.6316
    "J27" = 468.1875
    "K27" = 292.6316
    "L27" = 468.1875
    "M27" = -185.6316
    "N27" = -682.1875
    "H110" = 27398
    "J110" = 27398
    "L110" = 27398
    "N110" = -35578
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
